$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the 25x5 data grid (shared-string/text cells, matching the
# sharedStrings.xml produced upstream by the new calculation functions).
# The apostrophe prefix forces Excel to store the literal text rather than
# auto-converting the numeric-looking strings into numbers; we then restore
# the original cell style (Menlo 11 / black) so formatting is unaffected.
$c = $ws.Range("A1")
$c.Value = "'100.60351595169934"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B1")
$c.Value = "'102.14260620442913"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C1")
$c.Value = "'6.198468328403211"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D1")
$c.Value = "'136.30052724082003"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E1")
$c.Value = "'26.06016523583408"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A2")
$c.Value = "'1.0710223845459537"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B2")
$c.Value = "'2.2546358001145914"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C2")
$c.Value = "'0.010995968047632905"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D2")
$c.Value = "'0.07686091661226051"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E2")
$c.Value = "'0.008245096144436848"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A3")
$c.Value = "'984.7884985079821"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B3")
$c.Value = "'1737.1630289076818"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C3")
$c.Value = "'64.07778970717332"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D3")
$c.Value = "'2181.696925199133"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E3")
$c.Value = "'63.654699357652184"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A4")
$c.Value = "'-1393.4435139707673"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B4")
$c.Value = "'-2063.7872491995995"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C4")
$c.Value = "'-2043.5958308071545"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D4")
$c.Value = "'-2150.854119044336"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E4")
$c.Value = "'-1617.5949329424195"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A5")
$c.Value = "'4154680.3376139654"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B5")
$c.Value = "'310572316.36383325"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C5")
$c.Value = "'51.4692283885518"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D5")
$c.Value = "'34470.290793373584"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E5")
$c.Value = "'38.82610789940819"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A6")
$c.Value = "'46346.469752370685"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B6")
$c.Value = "'12.437377508369304"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C6")
$c.Value = "'254.3965458391769"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D6")
$c.Value = "'10.701316144738575"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E6")
$c.Value = "'2544.0785366615473"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A7")
$c.Value = "'3.6704183128172305"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B7")
$c.Value = "'4.136540012517814"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C7")
$c.Value = "'0.4065400158531448"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D7")
$c.Value = "'9.865781407943404"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E7")
$c.Value = "'1.5665397208456255"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A8")
$c.Value = "'12.330061772803079"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B8")
$c.Value = "'10.32952543019139"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C8")
$c.Value = "'9.104239817302151"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D8")
$c.Value = "'11.018995952407343"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E8")
$c.Value = "'9.176837973681687"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A9")
$c.Value = "'10994.53508334578"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B9")
$c.Value = "'11203.399306618961"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C9")
$c.Value = "'11321.644905422485"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D9")
$c.Value = "'11058.188469345123"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E9")
$c.Value = "'11461.055515913069"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A10")
$c.Value = "'1897070.1822955066"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B10")
$c.Value = "'187190152.0664467"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C10")
$c.Value = "'0.666666666666692"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D10")
$c.Value = "'5666.565412730209"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E10")
$c.Value = "'0.7620056364031144"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A11")
$c.Value = "'1.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B11")
$c.Value = "'1.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C11")
$c.Value = "'1.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D11")
$c.Value = "'1.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E11")
$c.Value = "'1.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A12")
$c.Value = "'-0.9999636452648891"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B12")
$c.Value = "'-0.9999552490545379"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C12")
$c.Value = "'-0.999999987108285"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D12")
$c.Value = "'-0.9999987843273923"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E12")
$c.Value = "'-0.9999999726141444"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A13")
$c.Value = "'-1.2990381001225255"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B13")
$c.Value = "'-1.2990381003993765"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C13")
$c.Value = "'-1.2990381056528009"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D13")
$c.Value = "'-1.299038105668708"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E13")
$c.Value = "'-1.2990381056766542"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A14")
$c.Value = "'0.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B14")
$c.Value = "'0.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C14")
$c.Value = "'0.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D14")
$c.Value = "'0.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E14")
$c.Value = "'0.0"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A15")
$c.Value = "'-0.9999730464833608"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B15")
$c.Value = "'-0.9999611685267041"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C15")
$c.Value = "'-0.9999999999999981"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D15")
$c.Value = "'-0.9999988959839906"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E15")
$c.Value = "'-0.9999999586088468"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A16")
$c.Value = "'8984.824996908215"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B16")
$c.Value = "'173282.6452646944"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C16")
$c.Value = "'428.82442286175404"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D16")
$c.Value = "'15139.384629004588"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E16")
$c.Value = "'1037.5463676988109"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A17")
$c.Value = "'3.9581557984026605"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B17")
$c.Value = "'8.661032775083093"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C17")
$c.Value = "'2.4825846680917506e-06"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D17")
$c.Value = "'9.20955333377097e-08"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E17")
$c.Value = "'3.0542986868152566e-20"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A18")
$c.Value = "'5584.81610528514"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B18")
$c.Value = "'365635.4753622076"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C18")
$c.Value = "'2.9709920721099626e-32"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D18")
$c.Value = "'4.3219772915880767e-39"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E18")
$c.Value = "'4.822732389775845e-98"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A19")
$c.Value = "'2.8375355834918587"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B19")
$c.Value = "'4.254722927801509"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C19")
$c.Value = "'1.1357056036115108e-11"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D19")
$c.Value = "'3.27752639842734e-13"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E19")
$c.Value = "'0.002263579591667563"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A20")
$c.Value = "'20.91518126671447"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B20")
$c.Value = "'20.251415027705345"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C20")
$c.Value = "'18.141908276339727"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D20")
$c.Value = "'20.861262671114034"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E20")
$c.Value = "'18.762232085279503"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A21")
$c.Value = "'17.591408090142448"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B21")
$c.Value = "'70.0877899763818"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C21")
$c.Value = "'1.0000000000116338"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D21")
$c.Value = "'1.000000000002697"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E21")
$c.Value = "'1.0000000000000584"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A22")
$c.Value = "'-0.007850064808855244"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B22")
$c.Value = "'-0.013333301672087029"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C22")
$c.Value = "'-0.04886067533344785"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D22")
$c.Value = "'-0.00812198530751035"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E22")
$c.Value = "'-0.029342795292674154"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A23")
$c.Value = "'-7.7416670858774195"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B23")
$c.Value = "'-18.038170483748498"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C23")
$c.Value = "'-14.085400524587198"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D23")
$c.Value = "'-7.819131323856326"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E23")
$c.Value = "'-11.252873440378412"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A24")
$c.Value = "'-6.153619167923483"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B24")
$c.Value = "'-6.216016845684947"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C24")
$c.Value = "'-6.133927383795085"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D24")
$c.Value = "'-6.30759295403667"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E24")
$c.Value = "'-6.307592954036664"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("A25")
$c.Value = "'-0.9999999980233191"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("B25")
$c.Value = "'-0.9999999975964072"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("C25")
$c.Value = "'-0.9999999999995917"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("D25")
$c.Value = "'-0.9999999999939"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0
$c = $ws.Range("E25")
$c.Value = "'-0.9999999999999988"
$c.Style = "Normal"
$c.Font.Name = "Menlo"
$c.Font.Size = 11
$c.Font.Color = 0

# Move the active selection to C8, matching the saved sheet view.
$ws.Range("C8").Select()
